$d = $word.ActiveDocument

$replacements = @(
    @("2024-12-02 Monday", "2024-12-03 Tuesday"),
    @("976÷4=", "199÷7="),
    @("869÷3=", "514÷5="),
    @("160÷4=", "486÷9="),
    @("753÷9=", "878÷4="),
    @("286÷4=", "275÷2="),
    @("326÷7=", "946÷8="),
    @("786÷6=", "623÷2="),
    @("693÷5=", "377÷7="),
    @("135÷5=", "810÷9="),
    @("488÷3=", "794÷3="),
    @("389÷4=", "337÷8="),
    @("617÷3=", "461÷6="),
    @("819÷4=", "201÷2="),
    @("304÷5=", "195÷3="),
    @("449÷6=", "558÷7="),
    @("281÷4=", "232÷9="),
    @("122÷4=", "898÷4="),
    @("835÷3=", "405÷7="),
    @("802÷3=", "545÷2="),
    @("232÷5=", "786÷8="),
    @("471÷3=", "887÷3="),
    @("313÷3=", "450÷4="),
    @("518÷7=", "822÷7="),
    @("288÷6=", "756÷9="),
    @("508÷8=", "147÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
